# Apply the edit described by the upstream diff:
#  - A50 ("Для автозаполнения ...") loses its trailing sentence about manual
#    registration (text shortened).
#  - A25 ("Нажмите на "Отправить контакт" ...") is shortened, dropping the
#    "или введите номер телефона" clause.
#  - Row 50's height shrinks (60pt -> 45pt) to match the now-shorter wrapped text.
#  - The sheet's scroll position / active selection moves up (it was parked at
#    A32 while viewing from row 19; now it's parked at A25 while viewing from
#    row 11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two long help-text cells with their shorter replacements.
# (A50 is written first so the new shared-string entries land in the same
# relative order the original workbook used.)
$ws.Range("A50").Value = 'Для автозаполнения основных данных из вашего аккаунта гелеграмм нажмите на кнопку "Отправить контакт". Я возьму имя, фамилию и телефон. '
$ws.Range("A25").Value = 'Нажмите на "Отправить контакт" .'

# The shorter text for A50 needs less vertical space.
$ws.Rows(50).RowHeight = 45

# Move the viewport / selection to match the edited workbook.
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A25").Select()
